$wb = $excel.ActiveWorkbook

# Grab the two sheets involved in this edit.
$validSheet = $wb.Worksheets.Item("ValidLoginTest")
$invalidSheet = $wb.Worksheets.Item("InvalidLoginTest")

# InvalidLoginTest loses focus; its selection moves to A4.
$invalidSheet.Range("A4").Select()

# ValidLoginTest becomes the active sheet.
$validSheet.Activate()

# Append a new "physician" test row to the valid-login test data.
$validSheet.Range("A4").Value = "physician"
$validSheet.Range("B4").Value = "physician"
$validSheet.Range("C4").Value = "OpenEMR"

# Leave the selection where Excel would land after entering the row via Tab/Enter.
$validSheet.Range("A5:XFD1048576").Select()
